$wb = $excel.ActiveWorkbook

# Rename the 'Codelists' sheet to 'Cells'
$ws = $wb.Worksheets.Item("Codelists")
$ws.Name = "Cells"

# Make 'Cells' the active sheet and select cell F13
$ws.Activate()
$ws.Range("F13").Select()
